$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -6.94354713543871
$ws.Range("C2").Value = 7.990168510511264
$ws.Range("D2").Value = 3.559831809881109
$ws.Range("B3").Value = -3.033464260293317
$ws.Range("C3").Value = 4.333352843923399
$ws.Range("D3").Value = -5.227846237893674
$ws.Range("B4").Value = -0.2234743598115374
$ws.Range("C4").Value = 4.093761345842939
$ws.Range("D4").Value = 1.620256385538821
$ws.Range("B5").Value = 1.260711594855279
$ws.Range("C5").Value = -1.223239834604506
$ws.Range("D5").Value = 8.070903554328556
$ws.Range("B6").Value = -5.160995246877953
$ws.Range("C6").Value = -3.632675769107285
$ws.Range("D6").Value = 0.1732976787794716
$ws.Range("B7").Value = -3.381558148626762
$ws.Range("C7").Value = 0.5958395557011942
$ws.Range("D7").Value = 0.7220832816449141
$ws.Range("B8").Value = -3.186190652963306
$ws.Range("C8").Value = -1.091206871444617
$ws.Range("D8").Value = -1.955526309528577
$ws.Range("B9").Value = 2.964647058339054
$ws.Range("C9").Value = 1.042552283805143
$ws.Range("D9").Value = 9.570679981139186
$ws.Range("B10").Value = -13.99768529350024
$ws.Range("C10").Value = -3.93685177458396
$ws.Range("D10").Value = -12.15959139070785
$ws.Range("B11").Value = -11.0410856605323
$ws.Range("C11").Value = 15.37783447774446
$ws.Range("D11").Value = -14.1120775080652
$ws.Range("B12").Value = -4.78117406122619
$ws.Range("C12").Value = 14.54741383364233
$ws.Range("D12").Value = -12.90242013598057
$ws.Range("B13").Value = -5.978844308965603
$ws.Range("C13").Value = 7.560548501577813
$ws.Range("D13").Value = -7.493641516982841
